# Auto update ESO load data
# Updates the hourly load values (columns D:AA) for existing rows 13-17
# and appends a brand new row 18 (2026-02-21) with a full day of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hourly columns D (hour 1) .. AA (hour 24)
$hourCols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")

# Updated hourly values for the already-present rows (13-17)
$rowUpdates = @{
    13 = @(4492,4268,4215,4156,4189,4332,4732,5158,5507,5674,5505,5510,5265,5263,5251,5319,5398,5597,6020,6039,5887,5670,5480,5148)
    14 = @(4821,4606,4511,4504,4513,4647,5054,5505,5854,5951,5849,5756,5643,5560,5444,5440,5477,5740,6322,6328,6141,5920,5739,5432)
    15 = @(4971,4756,4660,4653,4662,4796,5205,5657,6007,6104,5992,5890,5767,5676,5549,5545,5584,5862,6481,6487,6300,6078,5896,5588)
    16 = @(4968,4752,4656,4650,4658,4792,5201,5653,6003,6100,5988,5887,5763,5673,5546,5542,5581,5859,6476,6483,6296,6074,5892,5584)
    17 = @(4756,4541,4446,4439,4448,4581,4989,5439,5788,5885,5787,5698,5589,5510,5399,5395,5427,5659,6173,6179,5999,5786,5611,5316)
}

foreach ($r in $rowUpdates.Keys) {
    $values = $rowUpdates[$r]
    for ($i = 0; $i -lt $hourCols.Length; $i++) {
        $ws.Range($hourCols[$i] + $r).Value = $values[$i]
    }
}

# New row 18: 2026-02-21
$newRow = 18
$ws.Range("A" + $newRow).Value = 2026
$ws.Range("B" + $newRow).Value = 2
$ws.Range("C" + $newRow).Value = 21

$row18Values = @(4949,4750,4652,4612,4579,4659,4827,5228,5721,5717,5675,5736,5562,5403,5300,5362,5585,5814,6174,6106,5862,5711,5567,5428)
for ($i = 0; $i -lt $hourCols.Length; $i++) {
    $ws.Range($hourCols[$i] + $newRow).Value = $row18Values[$i]
}
